# agrego link a mi sitio en en y es
# The commit materializes the built-in "FollowedHyperlink" character style
# (normally latent until Word needs to persist it) into word/styles.xml,
# right after the last existing style ("UnresolvedMention"), with the same
# shape Word itself uses for the sibling "Hyperlink" style already present
# in this document.

$d = $word.ActiveDocument

# Create the new character style. wdStyleTypeCharacter = 2
$followed = $d.Styles.Add("FollowedHyperlink", 2)

# <w:basedOn w:val="DefaultParagraphFont"/>
$followed.BaseStyle = $d.Styles("DefaultParagraphFont")

# <w:uiPriority w:val="99"/>
$followed.Priority = 99

# <w:unhideWhenUsed/>
$followed.UnhideWhenUsed = $true

# <w:rPr><w:color w:val="954F72" w:themeColor="followedHyperlink"/><w:u w:val="single"/></w:rPr>
$followed.Font.TextColor.ObjectThemeColor = 11
$followed.Font.Underline = 1

Write-Output "Added FollowedHyperlink style"
